$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("B5","B8","B11","B12","B13","B14","B15","B16","B19","B27","B30","B31","B32","B33","B34")
foreach ($c in $cells) {
    $ws.Range($c).Value = "noir"
}

$cells = @("C2","C4","C10","C20")
foreach ($c in $cells) {
    $ws.Range($c).Value = "résultat postés ou publiés"
}

$cells = @("C3","C6","C7","C18","C22","C23","C24")
foreach ($c in $cells) {
    $ws.Range($c).Value = "résultat postés ou publiés dans les 36 mois"
}

$cells = @("C5","C8","C11","C12","C13","C14","C15","C16","C19","C27","C30","C31","C32","C33","C34")
foreach ($c in $cells) {
    $ws.Range($c).Value = "pas de résultat postés ni publiés"
}

$cells = @("C9","C17","C21","C25","C26","C28","C29")
foreach ($c in $cells) {
    $ws.Range($c).Value = "résultat postés ou publiés dans les 12 mois"
}
